# Refresh the cryptocurrency price/volume snapshot (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed "Price" figures are plain decimal numbers (e.g. "302.94").
# The sheet always stores that column as text, so force text formatting on
# those specific cells first to avoid Excel auto-converting them to numbers.
$textPriceRows = @(5, 6, 9, 10, 13, 15, 17, 20, 22, 23, 24, 29, 30, 33, 36, 39, 40, 42, 50, 51)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "42.967.75"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "2.359.88"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "302.94"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "95.24"
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  -3.23%  "

$ws.Range("D10").Value = "34.33"
$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("E11").Value = "  +2.03%  "

$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").Value = "18.49"
$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.726.23"
$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "6.69"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "2.368.57"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").Value = "42.956.16"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  +1.82%  "

$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "67.90"

$ws.Range("D23").Value = "235.12"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "2.19"
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("E28").Value = "  +15.71%  "

$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  +2.78%  "

$ws.Range("D30").Value = "32.24"
$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").Value = "17.51"
$ws.Range("E33").Value = "  +0.38%  "

$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("E35").Value = "  +6.69%  "

$ws.Range("D36").Value = "128.80"
$ws.Range("E36").Value = "  -7.84%  "

$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("D39").Value = "2.84"
$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("D40").Value = "2.26"
$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("E41").Value = "  -0.55%  "

$ws.Range("D42").Value = "20.73"
$ws.Range("E42").Value = "  -6.20%  "

$ws.Range("D43").Value = "1.927.54"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("E45").Value = "  +3.52%  "

$ws.Range("E46").Value = "  -9.51%  "

$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("D48").Value = "2.587.95"
$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("E49").Value = "  +3.11%  "

$ws.Range("D50").Value = "71.44"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "51.21"
$ws.Range("E51").Value = "  -2.51%  "
